# Update countries & provincias Spain
# Applies the refreshed COVID snapshot: new timestamp, updated totals for a
# few existing countries, and several countries whose case counts moved them
# past their neighbours in the (descending, by total cases) ranking -- which
# shifts who occupies each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 03:11"

# Estados Unidos (row 4) - updated totals, same rank
$ws.Range("B4").Value = 5655921
$ws.Range("C4").Value = 43946
$ws.Range("D4").Value = 3011098
$ws.Range("E4").Value = 2469756
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1351
$ws.Range("H4").Value = 175067

# India (row 6) - updated totals, same rank
$ws.Range("E6").Value = 676909
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1089
$ws.Range("H6").Value = 53014

# Libia jumps ahead of Guinea and Guayana Francesa (rows 92-94)
$ws.Range("A92").Value = "Libia"
$ws.Range("B92").Value = 9068
$ws.Range("C92").Value = 489
$ws.Range("D92").Value = 1003
$ws.Range("E92").Value = 7901
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 7
$ws.Range("H92").Value = 164

$ws.Range("A93").Value = "Guinea"
$ws.Range("B93").Value = 8715
$ws.Range("C93").Value = 95
$ws.Range("D93").Value = 7532
$ws.Range("E93").Value = 1131
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 52

$ws.Range("A94").Value = "Guayana Francesa"
$ws.Range("B94").Value = 8657
$ws.Range("C94").Value = 35
$ws.Range("D94").Value = 8054
$ws.Range("E94").Value = 550
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 53

# Bahamas jumps ahead of Malta, Republica de Chipre and Georgia (rows 145-148)
$ws.Range("A145").Value = "Bahamas"
$ws.Range("B145").Value = 1424
$ws.Range("C145").Value = 95
$ws.Range("D145").Value = 203
$ws.Range("E145").Value = 1201
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 20

$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 1423
$ws.Range("C146").Value = 48
$ws.Range("D146").Value = 766
$ws.Range("E146").Value = 648
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 9

$ws.Range("A147").Value = "Republica de Chipre"
$ws.Range("B147").Value = 1359
$ws.Range("C147").Value = 8
$ws.Range("D147").Value = 878
$ws.Range("E147").Value = 461
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 20

$ws.Range("A148").Value = "Georgia"
$ws.Range("B148").Value = 1351
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 1092
$ws.Range("E148").Value = 242
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 17

# Trinidad yTobago (row 166) - updated totals, same rank
$ws.Range("B166").Value = 629
$ws.Range("C166").Value = 41
$ws.Range("E166").Value = 477

# Barbados (row 188) - updated totals, same rank
$ws.Range("B188").Value = 153
$ws.Range("C188").Value = 1
$ws.Range("E188").Value = 24

# Santa Lucia jumps ahead of Timor Oriental (rows 202-203)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("B202").Value = 26
$ws.Range("C202").Value = 1
$ws.Range("D202").Value = 25
$ws.Range("E202").Value = 1

$ws.Range("A203").Value = "Timor Oriental"
$ws.Range("B203").Value = 25
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 24
$ws.Range("E203").Value = 1

# Islas Malvinas jumps ahead of Montserrat (rows 213-214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
